# Blue Ridge Community College Organizations - column rework
#  - Swap "Organization Name" (A) and "Categories" (B) columns
#  - Rename several header cells
#  - Add a new "Tiktok Link" column (M)
#  - Resize columns to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data in columns A and B (rows 2-18), keep header row separate ---
$colAVals = $ws.Range("A2:A18").Value2
$colBVals = $ws.Range("B2:B18").Value2
$ws.Range("A2:A18").Value2 = $colBVals
$ws.Range("B2:B18").Value2 = $colAVals

# --- Header row rewrite ---
$ws.Range("A1").Value2 = "Category"
$ws.Range("B1").Value2 = "Organization Name"
$ws.Range("C1").Value2 = "Organization Link"
$ws.Range("D1").Value2 = "Logo Link"
$ws.Range("E1").Value2 = "Description"
$ws.Range("F1").Value2 = "Email"
$ws.Range("G1").Value2 = "Phone Number"
$ws.Range("H1").Value2 = "Linkedin Link"
$ws.Range("I1").Value2 = "Instagram Link"
$ws.Range("J1").Value2 = "Facebook Link"
$ws.Range("K1").Value2 = "Twitter Link"
$ws.Range("L1").Value2 = "Youtube Link"

# New header cell: give it the same header formatting (bold, centered, bordered)
# as the rest of row 1 before filling in its text.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value2 = "Tiktok Link"

# --- New empty column M for data rows 2-18 ---
# A plain value-assignment of "" is indistinguishable from "no cell" to the
# engine, so the cell never gets materialised. Touching a formatting
# property with its already-default value forces the (still valueless)
# cell to be persisted, matching the source file's explicit empty cells.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 13).Font.Bold = $false
}

# --- Column widths (stored width = ColumnWidth + 5/6 for default Calibri 11) ---
$ws.Columns.Item(1).ColumnWidth = 17.166666666666668   # 18
$ws.Columns.Item(2).ColumnWidth = 49.166666666666664   # 50
$ws.Columns.Item(3).ColumnWidth = 49.166666666666664   # 50
$ws.Columns.Item(4).ColumnWidth = 49.166666666666664   # 50
$ws.Columns.Item(5).ColumnWidth = 49.166666666666664   # 50
$ws.Columns.Item(6).ColumnWidth = 22.166666666666668   # 23
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666   # 14
$ws.Columns.Item(8).ColumnWidth = 14.166666666666666   # 15
$ws.Columns.Item(9).ColumnWidth = 15.166666666666666   # 16
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666  # 15
$ws.Columns.Item(11).ColumnWidth = 13.166666666666666  # 14
$ws.Columns.Item(12).ColumnWidth = 13.166666666666666  # 14
$ws.Columns.Item(13).ColumnWidth = 12.166666666666666  # 13
